# Reorders the data rows of the active worksheet (component rows) so that
# each whole row (columns A:F) moves to a new row position. The values
# themselves are unchanged; only which row they sit in changes.
#
# Note: reading/writing a multi-cell Range via the ".Value" property is not
# reliable in this runtime, so ".Value2" is used instead (it correctly
# returns/accepts a 1-based 2D array for multi-cell ranges).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that participate in the reorder (two independent blocks).
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,18,19,20,21,22,23)

# Snapshot every row's current contents before any writes happen, so that
# moving data around does not clobber a row before it has been read.
$original = @{}
foreach ($r in $rows) {
    $original[$r] = $ws.Range("A" + $r + ":F" + $r).Value2
}

# Target row <- source row mapping (whole-row moves).
$mapping = @{
    2  = 12
    3  = 14
    4  = 13
    5  = 6
    6  = 2
    7  = 5
    8  = 10
    9  = 8
    10 = 4
    11 = 9
    12 = 7
    13 = 11
    14 = 3
    18 = 19
    19 = 21
    20 = 20
    21 = 18
    22 = 23
    23 = 22
}

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $ws.Range("A" + $target + ":F" + $target).Value2 = $original[$source]
}
